# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de locales have moved from "In Translation" to "Ready for handoff",
# refreshing the associated handoff timestamps, and widening the
# Status / handoff-datetime columns to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # "Overview"
$wsZhCn     = $wb.Worksheets.Item(2)   # "zh-cn"
$wsDeDe     = $wb.Worksheets.Item(3)   # "de-de"

# --- Overview sheet -------------------------------------------------------
# zh-cn / de-de status columns
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2017-02-21 10:48:08"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-02-21 10:47:51"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-02-21 10:48:08"

# --- Widen columns to fit the new "Ready for handoff" status text ----------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # de-de status column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33       # Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33       # Status column
